$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A, B, C, D, E ($null = leave blank)
$data = @(
    @("TEST", "BEAM", "UFF", "PLAN", "API PLAN"),
    @("beam 65 zip 58438 upfront fees 0 plan Augmented Standard Plan", "P", "P", "Augmented Standard Plan", "liberty12"),
    @("beam 65 zip 56740 upfront fees 0 plan Augmented Standard Plan", "P", "P", "Augmented Standard Plan", "liberty12"),
    @("beam 67 zip 49916 upfront fees 0 plan Augmented Standard Plan", "P", "P", "Augmented Standard Plan", "liberty12"),
    @("beam 79 zip 04763 upfront fees 0 plan Augmented Standard Plan", "P", "P", "Augmented Standard Plan", "liberty12"),
    @("beam 95 zip 33040 upfront fees 0 plan Augmented Standard Plan", "F", "F", "Augmented Standard Plan", "liberty12"),
    @("beam 301 zip 55806 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 302 zip 56440 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 303 zip 54701 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 304 zip 54124 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 305 zip 48661 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 309 zip 50126 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 310 zip 61032 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 311 zip 47963 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 312 zip 46777 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 314 zip 13730 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 315 zip 12208 upfront fees 150 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "liberty12-unlimited"),
    @("beam 317 zip 62702 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 319 zip 45414 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 322 zip 08724 upfront fees 150 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 325 zip 72401 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 327 zip 40744 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 328 zip 24121 upfront fees 150 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "liberty12-unlimited"),
    @("beam 330 zip 73505 upfront fees 150 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "liberty12-unlimited"),
    @("beam 332 zip 71603 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "gold30"),
    @("beam 333 zip 38901 upfront fees 0 plan VS2 Standard Offer", "P", "P", "VS2 Standard Offer", "gold30"),
    @("beam 334 zip 35758 upfront fees 0 plan Liberty 12/25/50/Unlimited", "P", "P", "Liberty 12/25/50/Unlimited", "liberty12-unlimited"),
    @("beam 336 zip 29204 upfront fees 0 plan Liberty 12/25/50/Unlimited", $null, $null, "Liberty 12/25/50/Unlimited", "gold30"),
)

# Row 1 (headers) is written left-to-right A,B,C,D,E.
# Data rows are written in the column order A, D, E, B, C, which matches how the
# original Katalon test-runner populated the sheet (test info first, then plan info,
# then the beam/uff pass-fail results last) -- this keeps the shared-string table order
# identical to the source workbook.
$headerOrder = @(1, 2, 3, 4, 5)
$dataOrder = @(1, 4, 5, 2, 3)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    $rowNum = $r + 1
    if ($rowNum -eq 1) {
        $order = $headerOrder
    } else {
        $order = $dataOrder
    }
    foreach ($col in $order) {
        $val = $row[$col - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}
